$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column W (23rd column), shifting the
# existing "eficiencia" .. "imagenes_carrusel" columns one position
# to the right (W->X, X->Y, ... AK->AL).
$ws.Range("W1").EntireColumn.Insert()

# New header for the inserted "codigo_de_producto" column.
$ws.Range("W1").Value = "codigo_de_producto"

# New data values for the two product rows.
$ws.Range("W2").Value = "SF-HEPA-H13-001"
$ws.Range("W3").Value = "SF-PRE-G4-001"
